$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the contents of a paragraph (excluding its end-of-paragraph
# mark) with a supplied OOXML <w:p> fragment. Using the sub-range that stops
# one character short of the paragraph mark preserves the paragraph's own
# identity/pPr while swapping out its runs for the new, more finely split
# run/proofErr structure.
# ---------------------------------------------------------------------------
function Set-ParagraphRuns {
    param([int]$Index, [string]$Xml)
    $para = $d.Paragraphs.Item($Index)
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    $rng = $d.Range($start, $end)
    $rng.InsertXML($Xml)
}

# ---------------------------------------------------------------------------
# Paragraph 6: "I squinted at this for a while and couldn't make anything of
# it. ... requests.get(url).text."
# Split the opening sentence around "couldn't" (proofErr gramStart/gramEnd)
# and split "requests.get(url).text" into "requests.get" (spellStart/End,
# not bold) + "(" + "url" (spellStart/End) + ").text" (still bold).
# ---------------------------------------------------------------------------
$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">I squinted at this for a while and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>couldn&#8217;t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> make anything of it. Then thought maybe &#8216;page source&#8217; in the hint means the webpage source. In my browser (Microsoft edge) I can right-click on a page and see the source. You can also do this with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>requests.get</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>).text</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@
Set-ParagraphRuns 6 $xml6

# ---------------------------------------------------------------------------
# Paragraph 7: "Anyway, looking at the source for the page, there's this:"
# Split around "there's" (proofErr gramStart/gramEnd).
# ---------------------------------------------------------------------------
$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Anyway, looking at the source for the page, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>there&#8217;s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> this:</w:t></w:r></w:p>
'@
Set-ParagraphRuns 7 $xml7

# ---------------------------------------------------------------------------
# Paragraph 14: "The only idea I could come up with was to count the
# frequency of all the characters in the mess. Printing these, the letters
# in 'equality' all have frequency 1. Plugging that into the url, it works.
# On to challenge 3."
# Split the second run around "url" (proofErr spellStart/spellEnd).
# ---------------------------------------------------------------------------
$xml14 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>The only idea I could come up with was to count the frequency of all the characters in the mess.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Printing these, the letters in &#8216;equality&#8217; all have frequency 1. Plugging that into the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>url</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, it works. On to challenge 3.</w:t></w:r></w:p>
'@
Set-ParagraphRuns 14 $xml14

# ---------------------------------------------------------------------------
# Remove the trailing empty paragraphs (previously 15-18): three blank
# Courier-New paragraphs plus the blank centred NormalWeb paragraph. The
# paragraph holding "On to challenge 3." now runs directly into the
# section properties.
# ---------------------------------------------------------------------------
$firstTrailing = $d.Paragraphs.Item(15)
$lastTrailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailingRange = $d.Range($firstTrailing.Range.Start, $lastTrailing.Range.End)
$trailingRange.Delete()
